$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.780.77'
$ws.Range('E2').Value = '  -1.76%  '

$ws.Range('D3').Value = '1.888.67'
$ws.Range('E3').Value = '  -1.77%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7733'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E8').Value = '  -4.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.25'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.33%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07203'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.65%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08087'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7646'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.75%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.478'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.58%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.894.68'
$ws.Range('E14').Value = '  -1.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.44%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.186'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.45%  '

$ws.Range('D17').Value = '29.795.41'
$ws.Range('E17').Value = '  -1.68%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.79%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007766'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '

$ws.Range('E21').Value = '  -0.01%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.147.76'
$ws.Range('E22').Value = '  -1.01%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.149'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.18%  '

$ws.Range('E24').Value = '  -0.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1582'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.26%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.437'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.04%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.40%  '

$ws.Range('E28').Value = '  -1.69%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.039'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.77%  '

$ws.Range('E30').Value = '  +4.92%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.547'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.460'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.37%  '

$ws.Range('E33').Value = '  -1.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05507'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.51%  '

$ws.Range('E35').Value = '  -3.84%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7508'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.19%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.002'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('E38').Value = '  -3.26%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01907'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.75%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.780'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.42%  '

$ws.Range('D41').Value = '1.155.43'
$ws.Range('E41').Value = '  +10.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.73%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4414'
$ws.Range('D43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.904'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.53%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8468'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.20%  '

$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('E47').Value = '  -2.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.85%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.979'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.56%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.447'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.88%  '

$ws.Range('E51').Value = '  -3.72%  '
